$d = $word.ActiveDocument

# --- 1. Expand the trailing whitespace after ${capHanhChinh} from 23 to 43 spaces ---
$oldCap = '${capHanhChinh}                       '
$newCap = '${capHanhChinh}                                           '
$d.Content.Find.Execute($oldCap, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newCap, 2)
